$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.160.58"
Set-TextValue "E2" "  -3.66%  "

Set-TextValue "D3" "3.296.36"
Set-TextValue "E3" "  -5.65%  "

Set-TextValue "E4" "  -0.07%  "

Set-TextValue "D5" "542.40"
Set-TextValue "E5" "  -2.35%  "

Set-TextValue "D6" "170.45"
Set-TextValue "E6" "  -4.61%  "

Set-TextValue "D7" "0.609"
Set-TextValue "E7" "  -4.44%  "

Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.13%  "

Set-TextValue "D9" "3.289.57"
Set-TextValue "E9" "  -5.76%  "

Set-TextValue "D10" "0.608"
Set-TextValue "E10" "  -3.67%  "

Set-TextValue "E11" "  -1.28%  "

Set-TextValue "D12" "52.35"
Set-TextValue "E12" "  -2.46%  "

Set-TextValue "D13" "0.0000262"
Set-TextValue "E13" "  -3.55%  "

Set-TextValue "D14" "8.82"
Set-TextValue "E14" "  -4.64%  "

Set-TextValue "D15" "3.814.53"
Set-TextValue "E15" "  -6.09%  "

Set-TextValue "E16" "  -2.41%  "

Set-TextValue "E17" "  -3.98%  "

Set-TextValue "D18" "3.287.73"
Set-TextValue "E18" "  -6.04%  "

Set-TextValue "D19" "11.59"
Set-TextValue "E19" "  -3.85%  "

Set-TextValue "D20" "63.078.66"
Set-TextValue "E20" "  -3.88%  "

Set-TextValue "E21" "  -3.47%  "

Set-TextValue "D22" "416.14"
Set-TextValue "E22" "  +0.59%  "

Set-TextValue "D23" "4.40"
Set-TextValue "E23" "  +7.07%  "

Set-TextValue "D24" "4.01"
Set-TextValue "E24" "  -0.94%  "

Set-TextValue "D25" "13.28"
Set-TextValue "E25" "  +4.57%  "

Set-TextValue "D26" "82.52"
Set-TextValue "E26" "  -4.03%  "

Set-TextValue "D27" "10.55"
Set-TextValue "E27" "  -2.34%  "

Set-TextValue "D28" "2.71"
Set-TextValue "E28" "  -4.92%  "

Set-TextValue "D29" "8.54"
Set-TextValue "E29" "  -5.49%  "

Set-TextValue "D30" "28.91"
Set-TextValue "E30" "  -4.48%  "

Set-TextValue "D31" "6.31"
Set-TextValue "E31" "  -2.55%  "

Set-TextValue "D32" "11.26"
Set-TextValue "E32" "  -3.35%  "

Set-TextValue "D33" "570.91"
Set-TextValue "E33" "  -6.30%  "

Set-TextValue "E34" "  -3.78%  "

Set-TextValue "D35" "57.81"
Set-TextValue "E35" "  -2.90%  "

Set-TextValue "E36" "  -0.11%  "

Set-TextValue "E37" "  -1.45%  "

Set-TextValue "E38" "  +4.30%  "

Set-TextValue "D39" "34.86"
Set-TextValue "E39" "  -6.18%  "

Set-TextValue "D40" "0.0₃0737"
Set-TextValue "E40" "  -6.49%  "

Set-TextValue "E41" "  -4.80%  "

Set-TextValue "D42" "3.098.43"
Set-TextValue "E42" "  -7.69%  "

Set-TextValue "D43" "0.999"
Set-TextValue "E43" "  -0.17%  "

Set-TextValue "D44" "2.74"
Set-TextValue "E44" "  -3.61%  "

Set-TextValue "E45" "  -2.10%  "

Set-TextValue "E48" "  -3.51%  "

Set-TextValue "E49" "  -5.51%  "

Set-TextValue "D50" "132.49"
Set-TextValue "E50" "  -3.72%  "

Set-TextValue "D51" "7.99"
Set-TextValue "E51" "  -5.34%  "

# Swap rows 46 and 47: Fetch.AI/VeChain reorder with new values
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D46" "0.0398"
Set-TextValue "E46" "  -3.96%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D47" "2.41"
Set-TextValue "E47" "  -4.59%  "
